$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $row, $col, $val) {
    # Force the value to be stored as text (many of these values look like
    # numbers, e.g. "1.000" or "0.00001008", but must remain literal strings),
    # then restore the "Normal" style so no stray formatting/style index is
    # left behind on the cell.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Rows 25 and 26: Cosmos/Stellar entries swapped places (with updated values)
Set-TextCell $ws 25 2 "Stellar"
Set-TextCell $ws 25 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws 25 4 "0.1386"
Set-TextCell $ws 25 5 "  +0.94%  "

Set-TextCell $ws 26 2 "Cosmos"
Set-TextCell $ws 26 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws 26 4 "8.509"
Set-TextCell $ws 26 5 "  +0.52%  "

# Price (column D) and Volume(1h) (column E) updates for remaining rows
$updates = @(
    @{Row=2;  D="29.048.63";     E="  -0.04%  "},
    @{Row=3;  D="1.833.73";      E="  +0.23%  "},
    @{Row=4;  D="0.9997";        E="  +0.08%  "},
    @{Row=5;  D="242.96";        E="  +0.63%  "},
    @{Row=6;  D="0.6266";        E="  -1.14%  "},
    @{Row=7;  D="1.000";         E="  +0.00%  "},
    @{Row=8;  D="0.07578";       E="  +3.20%  "},
    @{Row=9;  D="0.2926";        E="  -0.35%  "},
    @{Row=10; D="22.58";         E="  -1.30%  "},
    @{Row=11; D="0.07739";       E="  +0.75%  "},
    @{Row=12; D="1.836.96";      E="  +0.43%  "},
    @{Row=13; D="4.966";         E="  -0.48%  "},
    @{Row=14; D="0.6648";        E="  +0.19%  "},
    @{Row=15; D="0.00001008";    E="  +16.38%  "},
    @{Row=16; D="83.01";         E="  +1.26%  "},
    @{Row=17; D="6.070";         E="  +0.02%  "},
    @{Row=18; D="29.072.22";     E="  +0.55%  "},
    @{Row=19; D="226.80";        E="  +1.19%  "},
    @{Row=20; D="12.40";         E="  -0.15%  "},
    @{Row=21; D="1.002";         E="  +0.14%  "},
    @{Row=22; D="7.211";         E="  +1.19%  "},
    @{Row=23; D="1.001";         E="  +0.00%  "},
    @{Row=24; D="159.60";        E="  +1.03%  "},
    @{Row=27; D="17.94";         E="  +0.33%  "},
    @{Row=28; D="1.493";         E="  -0.63%  "},
    @{Row=29; D="4.098";         E="  +0.06%  "},
    @{Row=30; D="4.010";         E="  -0.39%  "},
    @{Row=31; D="1.193";         E="  -0.60%  "},
    @{Row=32; D="0.05251";       E="  -0.93%  "},
    @{Row=33; D="1.842";         E="  +0.57%  "},
    @{Row=34; D="0.7349";        E="  -0.55%  "},
    @{Row=35; D=$null;           E="  -1.50%  "},
    @{Row=36; D="2.690";         E="  +1.42%  "},
    @{Row=37; D="1.243.25";      E="  -3.82%  "},
    @{Row=38; D="2.763";         E="  +0.60%  "},
    @{Row=39; D="0.01783";       E="  -0.01%  "},
    @{Row=40; D="6.375";         E="  +1.45%  "},
    @{Row=41; D="0.8989";        E="  +0.17%  "},
    @{Row=42; D="1.000";         E="  +0.04%  "},
    @{Row=43; D="102.24";        E="  -0.45%  "},
    @{Row=44; D="1.985.19";      E="  +0.36%  "},
    @{Row=45; D=$null;           E="  +2.90%  "},
    @{Row=46; D="64.20";         E="  +0.15%  "},
    @{Row=47; D="0.5112";        E="  -0.45%  "},
    @{Row=48; D="0.4038";        E="  +1.24%  "},
    @{Row=49; D="8.882";         E="  +1.54%  "},
    @{Row=50; D="0.05764";       E="  -0.92%  "},
    @{Row=51; D="6.697";         E="  -0.04%  "}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        Set-TextCell $ws $r 4 $u.D
    }
    Set-TextCell $ws $r 5 $u.E
}
